# Commit: "Updated for 2 test cases"
#
# 1) LaunchPage/LaunchPageScript: fix a typo + append two new GooglePlayStore /
#    AppStore test-case rows.
# 2) Add two brand-new sheets (Login, LoginScript) mirroring the LaunchPage /
#    LaunchPageScript pair, for a "login" test flow.
# 3) SetupTests gains a second "login" suite toggle column, switched on.
# 4) The SetupTests tab becomes the active tab on save.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. LaunchPage (sheet3): rename the "Test"/"swiggy" scratch values to the
#    "DeliveryLocation"/"Jodhpur, Rajasthan" values actually used by the script.
# ---------------------------------------------------------------------------
$launchPage = $wb.Worksheets.Item("LaunchPage")
$launchPage.Range("A1").Value = "DeliveryLocation"
$launchPage.Range("A2").Value = "Jodhpur, Rajasthan"
$launchPage.Columns.Item(1).ColumnWidth = 17.15

# ---------------------------------------------------------------------------
# 2. LaunchPageScript (sheet4): fix the typo in A5, then append 4 new rows.
# ---------------------------------------------------------------------------
$launchPageScript = $wb.Worksheets.Item("LaunchPageScript")
$launchPageScript.Range("A5").Value = "Verify that text can be entered in the inputbox for delivery location."

$newRows = @(
    @(
        "Verify that the GooglePlayStore link is displayed on the homepage.",
        "GooglePlayStore link should be displayed on the homepage.",
        "GooglePlayStore link is displayed on the homepage.",
        "Problem finding the GooglePlayStore link."
    ),
    @(
        "Click on the GooglePlayStore link.",
        "A new tab should open with GooglePlayStore page for downloading Swiggy.",
        "A new tab opened with GooglePlayStore page for downloading Swiggy.",
        "Problem loading the GooglePlayStore page for Swiggy."
    ),
    @(
        "Verify the the AppStore link is displayed on the homepage.",
        "AppStore link should be displayed on the homepage.",
        "AppStore link is displayed on the homepage.",
        "Problem finding the AppStore link."
    ),
    @(
        "Click on the AppStore link.",
        "A new tab should open with AppStore page for downloading Swiggy.",
        "A new tab opened with AppStore page for downloading Swiggy.",
        "Problem loading the AppStore page for Swiggy."
    )
)

$r = 6
foreach ($row in $newRows) {
    $launchPageScript.Cells.Item($r, 1).Value = $row[0]
    $launchPageScript.Cells.Item($r, 2).Value = $row[1]
    $launchPageScript.Cells.Item($r, 3).Value = $row[2]
    $launchPageScript.Cells.Item($r, 4).Value = $row[3]
    $launchPageScript.Range($launchPageScript.Cells.Item($r, 1), $launchPageScript.Cells.Item($r, 4)).Style = $launchPageScript.Range("A3").Style
    $r = $r + 1
}

$launchPageScript.Range("A1:D3").Select()

# ---------------------------------------------------------------------------
# 3. New sheets: Login + LoginScript, appended after LaunchPageScript.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$login = $wb.Worksheets.Add($null, $lastSheet)
$login.Name = "Login"

$login.Range("A1").Value = "PhoneNo"
$login.Range("A2").Value = 7792016995
$login.Range("A1:A2").NumberFormat = "@"
$login.Range("B1").Value = "RandomOTP"
$login.Range("B2").Value = "123f"
$login.Columns.Item(1).AutoFit()
$login.Columns.Item(2).AutoFit()
$login.Range("A1:A1048576").Select()

$loginScript = $wb.Worksheets.Add($null, $login)
$loginScript.Name = "LoginScript"

$loginScript.Range("A1").Value = $launchPageScript.Range("A1").Value
$loginScript.Range("B1").Value = $launchPageScript.Range("B1").Value
$loginScript.Range("C1").Value = $launchPageScript.Range("C1").Value
$loginScript.Range("D1").Value = $launchPageScript.Range("D1").Value
$loginScript.Range("A1:D1").Style = $launchPageScript.Range("A1").Style

$loginScript.Range("A2").Value = $launchPageScript.Range("A2").Value
$loginScript.Range("B2").Value = $launchPageScript.Range("B2").Value
$loginScript.Range("C2").Value = $launchPageScript.Range("C2").Value
$loginScript.Range("D2").Value = $launchPageScript.Range("D2").Value
$loginScript.Range("B2").Style = $launchPageScript.Range("B2").Style
$loginScript.Rows.Item(2).RowHeight = 75

$loginRows = @(
    @(
        "Verify that the login link is present`n on the homepage.",
        "Login link should be present.",
        "Login link is present.",
        "Problem finding the login link on the homepage."
    ),
    @(
        "Click on the login link.",
        "A new form should appear for Login.",
        "Login form appeared.",
        "Problem loading the form for Login."
    ),
    @(
        "Enter the ten digit phone no in the Phone no field.",
        "Phone number should be entered.",
        "Phone number is entered.",
        "Problem entering phone number."
    ),
    @(
        "Click on login button.",
        "OTP field should be displayed.",
        "OTP field is displayed.",
        "Problem loading the OTP field."
    ),
    @(
        "Enter random number in OTP field and click Verify OTP.",
        "Enter Valid OTP warning should appear.",
        "Enter Valid OTP warning appears.",
        "Problem loading the warning."
    )
)

$r = 3
foreach ($row in $loginRows) {
    $loginScript.Cells.Item($r, 1).Value = $row[0]
    $loginScript.Cells.Item($r, 2).Value = $row[1]
    $loginScript.Cells.Item($r, 3).Value = $row[2]
    $loginScript.Cells.Item($r, 4).Value = $row[3]
    $loginScript.Range($loginScript.Cells.Item($r, 1), $loginScript.Cells.Item($r, 4)).Style = $launchPageScript.Range("A3").Style
    $loginScript.Rows.Item($r).RowHeight = 45
    $r = $r + 1
}

$loginScript.Columns.Item(1).AutoFit()
$loginScript.Columns.Item(2).AutoFit()
$loginScript.Columns.Item(3).AutoFit()
$loginScript.Columns.Item(4).AutoFit()
$loginScript.Range("A3").Select()

# ---------------------------------------------------------------------------
# 4. SetupTests (sheet2): add the "login" suite toggle column (enabled).
# ---------------------------------------------------------------------------
$setupTests = $wb.Worksheets.Item("SetupTests")
$setupTests.Range("B1").Value = "login"
$setupTests.Range("B2").Value = $true
$setupTests.Activate()
$setupTests.Range("A2").Select()

$wb.Saved = $false
